$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("analysis")

# Updated export timestamp (shared string used in A1)
$ws.Range("A1").Value = "25/03/2021 20:27"

# Updated occupancy values (recomputed export; overlaps now counted as int)
$ws.Range("C3").Value = 0.4503119129516114
$ws.Range("C4").Value = 0.5496880870483885

$ws.Range("D14").Value = 1

$ws.Range("C15").Value = 0.3365656838550645
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "0.0%"

$ws.Range("C16").Value = 0.3050197141048431
$ws.Range("D16").Value = 0.507602598195334

$ws.Range("C17").Value = 0.3564402734221735
$ws.Range("D17").Value = 0.478605129936724

$ws.Range("C18").Value = 0.3303662211294531
$ws.Range("D18").Value = 0.5077904797018222

$ws.Range("C19").Value = 0.06836831387729592
$ws.Range("C19").NumberFormat = "0.0%"
$ws.Range("D19").Value = 0.3561315541487868

$ws.Range("C20").Value = 0.4150792315697763
$ws.Range("D20").Value = 0.4896803191695611

$ws.Range("C21").Value = 0.490932171481582
$ws.Range("D21").Value = 0.426321057237913

$ws.Range("C22").Value = 0.7611659124538951
$ws.Range("D22").Value = 0.1447978060262054

$ws.Range("C23").Value = 0.7590097300531033
$ws.Range("D23").Value = 0.1071137301099014

$ws.Range("C24").Value = 0.7407743094437484
$ws.Range("D24").Value = 0.1530169980103788

$ws.Range("C25").Value = 0.8831022467414469
$ws.Range("D25").Value = 0.09766664629868996

$ws.Range("C26").Value = 0.9536766829963157
$ws.Range("D26").Value = 0.008770758953145222

$ws.Range("C27").Value = 0.9314503605489437
$ws.Range("D27").Value = 0
